# Translate the two column headers (B1 / C1) from Chinese to English.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Travel Man-Time"
$ws.Range("C1").Value = "Average Travel Distance"

# Widen column C so the longer English header fits (closest reachable
# width to the target 17.75 given this runtime's character-width rounding).
$ws.Columns("C").ColumnWidth = 17

# Update the embedded line chart: series names, axis titles and legend box.
$co = $ws.ChartObjects(1)
$chart = $co.Chart

# Series 1 plots column C ("Average Travel Distance"); series 2 plots
# column B ("Travel Man-Time").
$chart.SeriesCollection(1).Name = "Average Travel Distance"
$chart.SeriesCollection(2).Name = "Travel Man-Time"

# Category (X) axis title.
$chart.Axes(1, 1).AxisTitle.Text = "Month"

# Primary value (Y) axis title - shared with the "Average Travel Distance" series.
$chart.Axes(2, 1).AxisTitle.Text = "Average Travel Distance"

# Move/resize the legend box to accommodate the wider English series names.
$legend = $chart.Legend
$legend.Left = 0.70059182867578373
$legend.Width = 0.28492839950637383

# Move the active selection to C3, matching the saved cursor position.
$ws.Range("C3").Select() | Out-Null
